$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AssetVendor")
$ws.Columns.Item(2).Insert()
$ws.Range("C1:C5").Cut($ws.Range("B1"))
$ws.Range("C1:C5").ClearContents()
Write-Output "done"
